$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 12:45:56"
$ws.Cells.Item(3, 1).Value = "Total filas: 170"
$ws.Cells.Item(56, 1).Value = "07:38:39"
$ws.Cells.Item(56, 2).Value = "09:17"
$ws.Cells.Item(56, 3).Value = "14_ABASTO"
$ws.Cells.Item(56, 4).Value = 99
$ws.Cells.Item(56, 5).Value = "LP1912"
$ws.Cells.Item(57, 1).Value = "08:27:16"
$ws.Cells.Item(57, 2).Value = "09:17"
$ws.Cells.Item(57, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(57, 4).Value = 50
$ws.Cells.Item(57, 5).Value = "LP1912"
$ws.Cells.Item(111, 1).Value = "11:52:01"
$ws.Cells.Item(111, 2).Value = "11:54"
$ws.Cells.Item(111, 3).Value = "225_GOMEZ"
$ws.Cells.Item(111, 4).Value = 2
$ws.Cells.Item(111, 5).Value = "LP1912"
$ws.Cells.Item(113, 1).Value = "10:50:41"
$ws.Cells.Item(113, 2).Value = "11:54"
$ws.Cells.Item(113, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(113, 4).Value = 64
$ws.Cells.Item(113, 5).Value = "LP1912"
$ws.Cells.Item(133, 1).Value = "11:34:59"
$ws.Cells.Item(133, 2).Value = "12:35"
$ws.Cells.Item(133, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(133, 4).Value = 61
$ws.Cells.Item(133, 5).Value = "LP1912"
$ws.Cells.Item(134, 1).Value = "11:11:33"
$ws.Cells.Item(134, 2).Value = "12:35"
$ws.Cells.Item(134, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(134, 4).Value = 84
$ws.Cells.Item(134, 5).Value = "LP1912"
$ws.Cells.Item(153, 1).Value = "12:45:56"
$ws.Cells.Item(153, 2).Value = "13:18"
$ws.Cells.Item(153, 3).Value = "15_ABASTO"
$ws.Cells.Item(153, 4).Value = 33
$ws.Cells.Item(153, 5).Value = "LP1912"
$ws.Cells.Item(154, 1).Value = "12:45:56"
$ws.Cells.Item(154, 2).Value = "13:21"
$ws.Cells.Item(154, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(154, 4).Value = 36
$ws.Cells.Item(154, 5).Value = "LP1912"
$ws.Cells.Item(155, 1).Value = "11:54:18"
$ws.Cells.Item(155, 2).Value = "13:22"
$ws.Cells.Item(155, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(155, 4).Value = 88
$ws.Cells.Item(155, 5).Value = "LP1912"
$ws.Cells.Item(156, 1).Value = "11:34:59"
$ws.Cells.Item(156, 2).Value = "13:24"
$ws.Cells.Item(156, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(156, 4).Value = 110
$ws.Cells.Item(156, 5).Value = "LP1912"
$ws.Cells.Item(157, 1).Value = "11:47:17"
$ws.Cells.Item(157, 2).Value = "13:25"
$ws.Cells.Item(157, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(157, 4).Value = 98
$ws.Cells.Item(157, 5).Value = "LP1912"
$ws.Cells.Item(158, 1).Value = "12:11:52"
$ws.Cells.Item(158, 2).Value = "13:25"
$ws.Cells.Item(158, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(158, 4).Value = 74
$ws.Cells.Item(158, 5).Value = "LP1912"
$ws.Cells.Item(159, 1).Value = "11:34:59"
$ws.Cells.Item(159, 2).Value = "13:32"
$ws.Cells.Item(159, 3).Value = "215A_EL PATO"
$ws.Cells.Item(159, 4).Value = 118
$ws.Cells.Item(159, 5).Value = "LP1912"
$ws.Cells.Item(160, 1).Value = "12:11:52"
$ws.Cells.Item(160, 2).Value = "13:32"
$ws.Cells.Item(160, 3).Value = "14_ABASTO"
$ws.Cells.Item(160, 4).Value = 81
$ws.Cells.Item(160, 5).Value = "LP1912"
$ws.Cells.Item(161, 1).Value = "11:47:17"
$ws.Cells.Item(161, 2).Value = "13:33"
$ws.Cells.Item(161, 3).Value = "215A_EL PATO"
$ws.Cells.Item(161, 4).Value = 106
$ws.Cells.Item(161, 5).Value = "LP1912"
$ws.Cells.Item(162, 1).Value = "12:11:52"
$ws.Cells.Item(162, 2).Value = "13:46"
$ws.Cells.Item(162, 3).Value = "225_GOMEZ"
$ws.Cells.Item(162, 4).Value = 95
$ws.Cells.Item(162, 5).Value = "LP1912"
$ws.Cells.Item(163, 1).Value = "11:52:01"
$ws.Cells.Item(163, 2).Value = "13:47"
$ws.Cells.Item(163, 3).Value = "225_GOMEZ"
$ws.Cells.Item(163, 4).Value = 115
$ws.Cells.Item(163, 5).Value = "LP1912"
$ws.Cells.Item(164, 1).Value = "12:33:21"
$ws.Cells.Item(164, 2).Value = "13:54"
$ws.Cells.Item(164, 3).Value = "15_ABASTO"
$ws.Cells.Item(164, 4).Value = 81
$ws.Cells.Item(164, 5).Value = "LP1912"
$ws.Cells.Item(165, 1).Value = "12:11:52"
$ws.Cells.Item(165, 2).Value = "14:01"
$ws.Cells.Item(165, 3).Value = "10_OLMOS"
$ws.Cells.Item(165, 4).Value = 110
$ws.Cells.Item(165, 5).Value = "LP1912"
$ws.Cells.Item(166, 1).Value = "12:45:56"
$ws.Cells.Item(166, 2).Value = "14:01"
$ws.Cells.Item(166, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(166, 4).Value = 76
$ws.Cells.Item(166, 5).Value = "LP1912"
$ws.Cells.Item(167, 1).Value = "12:33:21"
$ws.Cells.Item(167, 2).Value = "14:02"
$ws.Cells.Item(167, 3).Value = "10_OLMOS"
$ws.Cells.Item(167, 4).Value = 89
$ws.Cells.Item(167, 5).Value = "LP1912"
$ws.Cells.Item(168, 1).Value = "12:45:56"
$ws.Cells.Item(168, 2).Value = "14:16"
$ws.Cells.Item(168, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(168, 4).Value = 91
$ws.Cells.Item(168, 5).Value = "LP1912"
$ws.Cells.Item(169, 1).Value = "12:33:21"
$ws.Cells.Item(169, 2).Value = "14:17"
$ws.Cells.Item(169, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(169, 4).Value = 104
$ws.Cells.Item(169, 5).Value = "LP1912"
$ws.Cells.Item(170, 1).Value = "12:33:21"
$ws.Cells.Item(170, 2).Value = "14:17"
$ws.Cells.Item(170, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(170, 4).Value = 104
$ws.Cells.Item(170, 5).Value = "LP1912"
$ws.Cells.Item(171, 1).Value = "12:45:56"
$ws.Cells.Item(171, 2).Value = "14:27"
$ws.Cells.Item(171, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(171, 4).Value = 102
$ws.Cells.Item(171, 5).Value = "LP1912"
$ws.Cells.Item(172, 1).Value = "12:45:56"
$ws.Cells.Item(172, 2).Value = "14:31"
$ws.Cells.Item(172, 3).Value = "14X44_ABASTO"
$ws.Cells.Item(172, 4).Value = 106
$ws.Cells.Item(172, 5).Value = "LP1912"
$ws.Cells.Item(173, 1).Value = "12:33:21"
$ws.Cells.Item(173, 2).Value = "14:32"
$ws.Cells.Item(173, 3).Value = "14X44_ABASTO"
$ws.Cells.Item(173, 4).Value = 119
$ws.Cells.Item(173, 5).Value = "LP1912"
$ws.Cells.Item(174, 1).Value = "12:45:56"
$ws.Cells.Item(174, 2).Value = "14:33"
$ws.Cells.Item(174, 3).Value = "215C_EL PATO"
$ws.Cells.Item(174, 4).Value = 108
$ws.Cells.Item(174, 5).Value = "LP1912"
$ws.Cells.Item(175, 1).Value = "12:45:56"
$ws.Cells.Item(175, 2).Value = "14:39"
$ws.Cells.Item(175, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(175, 4).Value = 114
$ws.Cells.Item(175, 5).Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 12:45:56"
$ws.Cells.Item(3, 1).Value = "Total filas: 27"
$ws.Cells.Item(32, 1).Value = "12:45:56"
$ws.Cells.Item(32, 2).Value = "14:33"
$ws.Cells.Item(32, 3).Value = "215C_EL PATO"
$ws.Cells.Item(32, 4).Value = 108
$ws.Cells.Item(32, 5).Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 12:45:56"
